$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. H column "score" cells: rows that get a numeric value, filled (theme9) + centered ---
$valueRows = @(3,5,7,9,11,15,17,19,21,23,25,38)
$valueMap = @{3=1;5=2;7=3;9=4;11=5;15=6;17=7;19=8;21=8;23=10;25=11;38=12}

# Union range of all the "scored" H cells, including the still-blank H27, so they
# all share one centred+filled style (matches the workbook's fillId=6 look).
$scoreRange = $ws.Range("H3,H5,H7,H9,H11,H15,H17,H19,H21,H23,H25,H27,H38")
$scoreRange.HorizontalAlignment = -4108
$scoreRange.VerticalAlignment = -4108
$scoreRange.Interior.ThemeColor = 9
$scoreRange.Interior.TintAndShade = 0

foreach ($r in $valueRows) {
    $ws.Range("H$r").Value = $valueMap[$r]
}
# H27 stays empty (no v element) but keeps the same formatting as the other score cells.

# --- 2. Separator rows (13, 29, 34): blank but filled (theme5) + centered, like the row's own style ---
$sepRange = $ws.Range("H13,H29,H34")
$sepRange.HorizontalAlignment = -4108
$sepRange.VerticalAlignment = -4108
$sepRange.Interior.ThemeColor = 5
$sepRange.Interior.TintAndShade = 0.79998168889431442

# --- 3. Plain filler H cells that simply extend the row to column H, matching the
#        centred / unfilled look already used by columns B:G on those rows ---
$plainRows = @(4,6,8,10,12,14,16,18,20,22,24,26,28,30,31,32,33,35,36,37,39,40,41,42,43)
$plainAddrs = ($plainRows | ForEach-Object { "H$_" }) -join ","
$plainRange = $ws.Range($plainAddrs)
$plainRange.HorizontalAlignment = -4108
$plainRange.VerticalAlignment = -4108
$plainRange.Interior.Pattern = -4142

# --- 4. Sheet view: scroll position + current selection moved further down the list ---
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("H38").Select()
